$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D (Price) and E (Volume) columns are treated as plain text so that
# values like "1.00", "0.000270" etc. keep their exact original formatting
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.953.72"
$ws.Range("E2").Value = "  -7.16%  "

$ws.Range("D3").Value = "3.689.31"
$ws.Range("E3").Value = "  -7.63%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "578.94"
$ws.Range("E5").Value = "  -5.18%  "

$ws.Range("E6").Value = "  +5.21%  "

$ws.Range("D7").Value = "3.689.46"
$ws.Range("E7").Value = "  -7.39%  "

$ws.Range("D8").Value = "0.632"
$ws.Range("E8").Value = "  -7.50%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("D10").Value = "0.717"
$ws.Range("E10").Value = "  -5.24%  "

$ws.Range("D11").Value = "0.166"
$ws.Range("E11").Value = "  -3.72%  "

$ws.Range("D12").Value = "53.02"
$ws.Range("E12").Value = "  -6.85%  "

$ws.Range("E13").Value = "  -11.38%  "

$ws.Range("D14").Value = "10.74"
$ws.Range("E14").Value = "  -3.36%  "

$ws.Range("D15").Value = "4.276.54"
$ws.Range("E15").Value = "  -7.61%  "

$ws.Range("D16").Value = "3.711.58"
$ws.Range("E16").Value = "  -6.97%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "19.42"
$ws.Range("E17").Value = "  -5.68%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.127"
$ws.Range("E18").Value = "  -3.15%  "

$ws.Range("D19").Value = "13.09"
$ws.Range("E19").Value = "  -8.09%  "

$ws.Range("E20").Value = "  -9.45%  "

$ws.Range("D21").Value = "67.801.51"
$ws.Range("E21").Value = "  -7.20%  "

$ws.Range("D22").Value = "410.39"
$ws.Range("E22").Value = "  -6.72%  "

$ws.Range("D23").Value = "4.56"
$ws.Range("E23").Value = "  -7.10%  "

$ws.Range("D24").Value = "89.26"
$ws.Range("E24").Value = "  -7.12%  "

$ws.Range("D25").Value = "3.07"
$ws.Range("E25").Value = "  -9.30%  "

$ws.Range("D26").Value = "12.82"
$ws.Range("E26").Value = "  -9.88%  "

$ws.Range("D27").Value = "10.87"
$ws.Range("E27").Value = "  -2.10%  "

$ws.Range("D28").Value = "3.80"
$ws.Range("E28").Value = "  -7.44%  "

$ws.Range("D29").Value = "5.91"
$ws.Range("E29").Value = "  -1.07%  "

$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "8.32"
$ws.Range("E30").Value = "  +8.54%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "9.63"
$ws.Range("E31").Value = "  -8.88%  "

$ws.Range("D32").Value = "33.02"
$ws.Range("E32").Value = "  -8.61%  "

$ws.Range("D33").Value = "12.77"
$ws.Range("E33").Value = "  -7.07%  "

$ws.Range("D34").Value = "69.23"
$ws.Range("E34").Value = "  -4.11%  "

$ws.Range("D35").Value = "45.06"
$ws.Range("E35").Value = "  -6.20%  "

$ws.Range("D36").Value = "0.119"
$ws.Range("E36").Value = "  -8.63%  "

$ws.Range("D37").Value = "0.0₃0933"
$ws.Range("E37").Value = "  -9.00%  "

$ws.Range("D38").Value = "584.05"
$ws.Range("E38").Value = "  -8.00%  "

$ws.Range("D39").Value = "0.405"
$ws.Range("E39").Value = "  -6.17%  "

$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("D42").Value = "3.23"
$ws.Range("E42").Value = "  +12.48%  "

$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  -7.28%  "

$ws.Range("D44").Value = "3.07"
$ws.Range("E44").Value = "  -10.49%  "

$ws.Range("D45").Value = "0.0443"
$ws.Range("E45").Value = "  -8.66%  "

$ws.Range("D46").Value = "2.62"
$ws.Range("E46").Value = "  +0.36%  "

$ws.Range("D47").Value = "9.48"
$ws.Range("E47").Value = "  -14.18%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "3.26"
$ws.Range("E48").Value = "  -5.32%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "27.41"
$ws.Range("E49").Value = "  +8.48%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.136"
$ws.Range("E50").Value = "  -8.78%  "

$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "0.000270"
$ws.Range("E51").Value = "  -3.50%  "
